# This workbook logs daily grape ("Uva") price records for "Terminal La
# Palmera de La Serena". A new record (dated 2022-03-18 / serial 44638) is
# being added to the weekly log. It belongs right above the existing row 43
# entry, so that row and everything below it (through row 99) shifts down by
# one, and the used range grows from A1:T99 to A1:T100.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 43 - this shifts rows 43:99 down to 44:100 and
# expands the sheet's used range accordingly (it also carries the column D
# date-cell formatting down into the new row, same as Excel's native
# "Insert Cut Cells" / "Insert Sheet Rows" behavior).
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new price record.
$ws.Range("A43").Value = 8
$ws.Range("B43").Value = "Terminal La Palmera de La Serena"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 44638
$ws.Range("E43").Value = 4
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100109
$ws.Range("H43").Value = "Uva"
$ws.Range("I43").Value = 100109001
$ws.Range("J43").Value = "Uva"
$ws.Range("K43").Value = "Red Globe"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 400
$ws.Range("N43").Value = 9500
$ws.Range("O43").Value = 10000
$ws.Range("P43").Value = 9750
$ws.Range("Q43").Value = "`$/bandeja 18 kilos"
$ws.Range("R43").Value = "Provincia del Elquí"
$ws.Range("S43").Value = 542
$ws.Range("T43").Value = 18
